$wb = $excel.ActiveWorkbook

# Row 33 on ALC (hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 483.5
$ws.Range("I33").Value = 472.8
$ws.Range("J33").Value = 515.6
$ws.Range("K33").Value = 472.8
$ws.Range("L33").Value = 515.6
$ws.Range("M33").Value = -243.8
$ws.Range("N33").Value = -973.6

# Row 137 on ALC (hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1269.1082
$ws.Range("I137").Value = 1119.1666
$ws.Range("J137").Value = 1545.9231
$ws.Range("K137").Value = 3357.4998
$ws.Range("L137").Value = 4637.7693
$ws.Range("M137").Value = -807.4998000000001
$ws.Range("N137").Value = -9737.7693

# Row 6 on ARM (hunk 2)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = $null

# Row 32 on ARM (hunk 3)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5329.01
$ws.Range("I32").Value = 3953.8408
$ws.Range("J32").Value = 15413.583
$ws.Range("K32").Value = 3953.8408
$ws.Range("L32").Value = 15413.583
$ws.Range("M32").Value = -3666.8408
$ws.Range("N32").Value = -15987.583

# Row 61 on ARM (hunk 4)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 66667988
$ws.Range("I61").Value = 76924136
$ws.Range("K61").Value = 76924136
$ws.Range("M61").Value = -76923924

# Row 74 on ARM (hunk 5)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3884
$ws.Range("I74").Value = 3794.75
$ws.Range("J74").Value = 3923.6667
$ws.Range("K74").Value = 3794.75
$ws.Range("L74").Value = 3923.6667
$ws.Range("M74").Value = -2920.75
$ws.Range("N74").Value = -5671.6667

# Row 77 on ARM (hunk 6)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3884
$ws.Range("I77").Value = 3794.75
$ws.Range("J77").Value = 3923.6667
$ws.Range("K77").Value = 18973.75
$ws.Range("L77").Value = 19618.3335
$ws.Range("M77").Value = -14605.75
$ws.Range("N77").Value = -28354.3335

# Row 110 on ARM (hunk 7)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 2164.7
$ws.Range("J110").Value = 3453.25
$ws.Range("L110").Value = 3453.25
$ws.Range("N110").Value = -7543.25

# Row 122 on ARM (hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2262.3684
$ws.Range("I122").Value = 1944.5294
$ws.Range("J122").Value = 4964
$ws.Range("K122").Value = 5833.5882
$ws.Range("L122").Value = 14892
$ws.Range("M122").Value = -3383.5882
$ws.Range("N122").Value = -19792

# Row 136 on ARM (hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 66667988
$ws.Range("I136").Value = 76924136
$ws.Range("K136").Value = 230772408
$ws.Range("M136").Value = -230769858

# Row 20 on BSM (hunk 10)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1024.138
$ws.Range("I20").Value = 1015.4231
$ws.Range("J20").Value = 1099.6666
$ws.Range("K20").Value = 1015.4231
$ws.Range("L20").Value = 1099.6666
$ws.Range("M20").Value = -768.4231
$ws.Range("N20").Value = -1593.6666

# Row 92 on BSM (hunk 11)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 14500
$ws.Range("J92").Value = 14500
$ws.Range("L92").Value = 14500
$ws.Range("N92").Value = -19492

# Row 99 on BSM (hunk 12)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 83334450
$ws.Range("I99").Value = 100001040
$ws.Range("K99").Value = 100001040
$ws.Range("M99").Value = -99999542

# Row 135 on BSM (hunk 13)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 34563.332
$ws.Range("J135").Value = 34563.332
$ws.Range("L135").Value = 34563.332
$ws.Range("N135").Value = -44703.332

# Row 31 on CRP (hunk 14)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1226.0597
$ws.Range("I31").Value = 1184.5186
$ws.Range("J31").Value = 1398.6154
$ws.Range("K31").Value = 1184.5186
$ws.Range("L31").Value = 1398.6154
$ws.Range("M31").Value = -889.5186000000001
$ws.Range("N31").Value = -1988.6154

# Row 34 on CRP (hunk 15)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1226.0597
$ws.Range("I34").Value = 1184.5186
$ws.Range("J34").Value = 1398.6154
$ws.Range("K34").Value = 1184.5186
$ws.Range("L34").Value = 1398.6154
$ws.Range("M34").Value = -982.5186000000001
$ws.Range("N34").Value = -1802.6154

# Row 56 on CRP (hunk 16)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 10000
$ws.Range("I56").Value = 10000
$ws.Range("K56").Value = 10000
$ws.Range("M56").Value = -9155

# Row 58 on CRP (hunk 17)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6731.5454
$ws.Range("I58").Value = 1162
$ws.Range("K58").Value = 1162
$ws.Range("M58").Value = -959

# Row 107 on CRP (hunk 18)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 815.2632
$ws.Range("I107").Value = 427.85715
$ws.Range("J107").Value = 1900
$ws.Range("K107").Value = 427.85715
$ws.Range("L107").Value = 1900
$ws.Range("M107").Value = 1492.14285
$ws.Range("N107").Value = -5740

# Row 134 on CRP (hunk 19)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 23811218
$ws.Range("I134").Value = 1635.5
$ws.Range("J134").Value = 100001880
$ws.Range("K134").Value = 4906.5
$ws.Range("L134").Value = 300005640
$ws.Range("M134").Value = -2371.5
$ws.Range("N134").Value = -300010710

# Row 136 on CRP (hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6731.5454
$ws.Range("I136").Value = 1162
$ws.Range("K136").Value = 3486
$ws.Range("M136").Value = -936

# Row 8 on CUL (hunk 21)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 124.9
$ws.Range("I8").Value = 124.9
$ws.Range("K8").Value = 374.7
$ws.Range("M8").Value = -235.7

# Row 23 on CUL (hunk 22)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 286.30768
$ws.Range("I23").Value = 68.2
$ws.Range("J23").Value = 422.625
$ws.Range("K23").Value = 204.6
$ws.Range("L23").Value = 1267.875
$ws.Range("M23").Value = 30.39999999999998
$ws.Range("N23").Value = -1737.875

# Row 131 on CUL (hunk 23)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 24393544
$ws.Range("J131").Value = 4238.161
$ws.Range("L131").Value = 12714.483
$ws.Range("N131").Value = -22794.483

# Row 7 on GSM (hunk 24)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5571428.5
$ws.Range("I7").Value = 5615384.5
$ws.Range("J7").Value = 5000000
$ws.Range("K7").Value = 5615384.5
$ws.Range("L7").Value = 5000000
$ws.Range("M7").Value = -5615272.5
$ws.Range("N7").Value = -5000224

# Row 8 on GSM (hunk 25)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H8").Value = 5571428.5
$ws.Range("I8").Value = 5615384.5
$ws.Range("J8").Value = 5000000
$ws.Range("K8").Value = 5615384.5
$ws.Range("L8").Value = 5000000
$ws.Range("M8").Value = -5615245.5
$ws.Range("N8").Value = -5000278

# Row 11 on GSM (hunk 26)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 7764706
$ws.Range("I11").Value = 7714286
$ws.Range("J11").Value = 8000000
$ws.Range("K11").Value = 7714286
$ws.Range("L11").Value = 8000000
$ws.Range("M11").Value = -7714147
$ws.Range("N11").Value = -8000278

# Row 41 on GSM (hunk 27)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1100
$ws.Range("I41").Value = 1100
$ws.Range("K41").Value = 1100
$ws.Range("M41").Value = -745

# Row 113 on GSM (hunk 28)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1831.1538
$ws.Range("I113").Value = 1837.9166
$ws.Range("K113").Value = 1837.9166
$ws.Range("M113").Value = 332.0834

# Row 100 on LTW (hunk 29)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1173.3334
$ws.Range("I100").Value = 942.8570999999999
$ws.Range("J100").Value = 1496
$ws.Range("K100").Value = 942.8570999999999
$ws.Range("L100").Value = 1496
$ws.Range("M100").Value = -401.8570999999999
$ws.Range("N100").Value = -2578

# Row 132 on LTW (hunk 30)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3659.7334
$ws.Range("J132").Value = 3119.8
$ws.Range("L132").Value = 9359.400000000001
$ws.Range("N132").Value = -14419.4

# Row 136 on LTW (hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2055.8948
$ws.Range("I136").Value = 1923.1666
$ws.Range("K136").Value = 5769.4998
$ws.Range("M136").Value = -3219.4998

# Row 132 on WVR (hunk 32)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1958.7894
$ws.Range("I132").Value = 1804.8064
$ws.Range("K132").Value = 5414.4192
$ws.Range("M132").Value = -2884.4192

# Row 136 on WVR (hunk 33)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1474.8096
$ws.Range("I136").Value = 1261.4
$ws.Range("J136").Value = 1668.8182
$ws.Range("K136").Value = 3784.2
$ws.Range("L136").Value = 5006.4546
$ws.Range("M136").Value = -1234.2
$ws.Range("N136").Value = -10106.4546
